# Rever_DailyTrack_BALRAJ_2022.xlsx update
# - Adds task #9 (rows 22-23) to the MAR-22 sheet:
#     Row 22: No=9, Date=2022-03-11, Application=RPA GSS,
#             Task="1. Service Order Pending task is work in progress , whereas, downloading and uploading testing is going",
#             % of completion=80%, Status=WIP
#     Row 23: Application=RPA RLOGIC,
#             Task="1. General Ledger reports has been triggered today for all three centers.",
#             % of completion=100%, Status=Completed
# - Moves the sheet view down to keep the newly-entered rows in view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MAR-22")

# --- Row 22 ---------------------------------------------------------------
$ws.Cells.Item(22, 1).Value = 9

# Copy the date format (numFmtId 14, m/d/yyyy) from an existing date cell
# so the same style slot is reused instead of a new one being created.
$ws.Cells.Item(2, 2).Copy()
$ws.Cells.Item(22, 2).PasteSpecial(-4122)
$ws.Cells.Item(22, 2).Value = 44631

$ws.Cells.Item(22, 3).Value = "RPA GSS"
$ws.Cells.Item(22, 4).Value = "1. Service Order Pending task is work in progress , whereas, downloading and uploading testing is going"

# Copy the percentage format (numFmtId 9, 0%) from an existing cell.
$ws.Cells.Item(2, 5).Copy()
$ws.Cells.Item(22, 5).PasteSpecial(-4122)
$ws.Cells.Item(22, 5).Value = 0.8

$ws.Cells.Item(22, 6).Value = "WIP"

# --- Row 23 ---------------------------------------------------------------
$ws.Cells.Item(23, 3).Value = "RPA RLOGIC"
$ws.Cells.Item(23, 4).Value = "1. General Ledger reports has been triggered today for all three centers."

$ws.Cells.Item(2, 5).Copy()
$ws.Cells.Item(23, 5).PasteSpecial(-4122)
$ws.Cells.Item(23, 5).Value = 1

$ws.Cells.Item(23, 6).Value = "Completed"

$excel.CutCopyMode = $false

# --- Scroll / selection update ---------------------------------------------
$ws.Range("D28").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 13
